# Add a new "50 scale FOCUS set" print request to the "August 2018" sheet
# (already the active sheet/tab), appending 4 new rows below the existing
# data (rows 2-21) as rows 22-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: FOCUS Box bot
$ws.Range("A22").Value = "16-08-2018"
$ws.Range("B22").Value = "16-08-2018"
$ws.Range("C22").Value = "FOCUS Box bot"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "PLA"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 20
$ws.Range("H22").Value = 0.2

# Row 23: 50 scale LAX
$ws.Range("A23").Value = "16-08-2018"
$ws.Range("C23").Value = "50 scale LAX"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = "Polylite"
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 0.2

# Row 24: 50 scale SAX
$ws.Range("A24").Value = "16-08-2018"
$ws.Range("C24").Value = "50 scale SAX"
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = "Polylite"
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 20
$ws.Range("H24").Value = 0.2

# Row 25: FOCUS Box tops
$ws.Range("A25").Value = "16-08-2018"
$ws.Range("C25").Value = "FOCUS Box tops"
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = "PLA"
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 0.2

# Match the author's final selection/cursor position
[void]$ws.Range("A25").Select()
